$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.258.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.21%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.587.25'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.37%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '508.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.13'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.19%  '
$ws.Range("E7").Value = '  +0.28%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.93%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.593.73'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.72%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.68'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +6.69%  '
$ws.Range("E11").Value = '  -1.22%  '
$ws.Range("E12").Value = '  -0.41%  '
$ws.Range("E13").Value = '  +1.87%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.040.05'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.49%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.225.11'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.46'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.46%  '
$ws.Range("E17").Value = '  +0.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.592.24'
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '353.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.44%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.87%  '
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.28'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.419'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.32%  '
$ws.Range("E26").Value = '  +0.24%  '
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0832'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.50%  '
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '19.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.78%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '151.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.45%  '
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.70'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("E35").Value = '  -1.05%  '
$ws.Range("E36").Value = '  -2.88%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.864'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.84%  '
$ws.Range("E38").Value = '  -2.30%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.76%  '
$ws.Range("E40").Value = '  -0.23%  '
$ws.Range("E41").Value = '  -2.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '294.07'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.67%  '
$ws.Range("E43").Value = '  -0.68%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.615'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.52%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.997'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.58%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0550'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.88%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.57'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.43%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.77'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.84%  '
$ws.Range("E49").Value = '  -1.67%  '
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.986.02'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.39%  '
